# Adds the "Sunday, Jan 15" departures (rows 315-346) to the Main Data sheet,
# continuing the existing NUMBER/DATE/TIME/FLIGHT/TO/SHORT/AIRLINE/MODEL/
# AIRCFAT ID/STATUS/DIFFERENCE table (columns A-J and L; K stays blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, NUMBER(A), DATE(B), TIME(C), FLIGHT(D), TO(E), SHORT(F),
#             AIRLINE(G), MODEL(H), AIRCFAT ID(I), STATUS(J), DIFFERENCE(L)
$newRows = @(
    @(315, 314, 'Sunday, Jan 15', '5:35 AM', 'LO3828', 'Warsaw', '(WAW)', 'LOT ', 'E190', '(SP-LMG)', '5:44 AM', '0 hours, 9 minutes'),
    @(316, 315, 'Sunday, Jan 15', '6:00 AM', 'KL1916', 'Amsterdam', '(AMS)', 'KLM ', 'E295', '(PH-NXN)', '6:22 AM', '0 hours, 22 minutes'),
    @(317, 316, 'Sunday, Jan 15', '6:00 AM', 'W61601', 'London', '(LTN)', 'Wizz Air ', 'A320', '(HA-LYH)', '6:09 AM', '0 hours, 9 minutes'),
    @(318, 317, 'Sunday, Jan 15', '6:05 AM', 'FR6125', 'Malta', '(MLA)', 'Ryanair ', 'B738', '(SP-RSL)', '6:11 AM', '0 hours, 6 minutes'),
    @(319, 318, 'Sunday, Jan 15', '6:05 AM', 'W61743', 'Oslo', '(TRF)', 'Wizz Air ', 'A320', '(HA-LYO)', '6:16 AM', '0 hours, 11 minutes'),
    @(320, 319, 'Sunday, Jan 15', '6:10 AM', 'FR6099', 'Stockholm', '(ARN)', 'Ryanair ', 'B738', '(SP-RKM)', '6:19 AM', '0 hours, 9 minutes'),
    @(321, 320, 'Sunday, Jan 15', '6:15 AM', 'W61739', 'Malmo', '(MMX)', 'Wizz Air ', 'A321', '(HA-LTB)', '6:21 AM', '0 hours, 6 minutes'),
    @(322, 321, 'Sunday, Jan 15', '6:25 AM', 'W61641', 'Eindhoven', '(EIN)', 'Wizz Air ', 'A21N', '(9H-WBU)', '6:36 AM', '0 hours, 11 minutes'),
    @(323, 322, 'Sunday, Jan 15', '6:40 AM', 'FR2375', 'London', '(STN)', 'Ryanair ', 'B738', '(SP-RKQ)', '6:54 AM', '0 hours, 14 minutes'),
    @(324, 323, 'Sunday, Jan 15', '7:05 AM', 'FR3283', 'Leeds', '(LBA)', 'Ryanair ', 'B738', '(SP-RSO)', '7:10 AM', '0 hours, 5 minutes'),
    @(325, 324, 'Sunday, Jan 15', '7:05 AM', 'FR3686', 'Billund', '(BLL)', 'Ryanair ', 'B738', '(SP-RSW)', '7:06 AM', '0 hours, 1 minutes'),
    @(326, 325, 'Sunday, Jan 15', '7:10 AM', 'W61661', 'Dortmund', '(DTM)', 'Wizz Air ', 'A320', '(HA-LYM)', '7:08 AM', '0 hours, -2 minutes'),
    @(327, 326, 'Sunday, Jan 15', '7:15 AM', 'W61745', 'Bergen', '(BGO)', 'Wizz Air ', 'A321', '(HA-LXL)', '8:08 AM', '0 hours, 53 minutes'),
    @(328, 327, 'Sunday, Jan 15', '7:50 AM', 'LH1647', 'Munich', '(MUC)', 'Lufthansa ', 'CRJ9', '(D-ACKH)', '7:52 AM', '0 hours, 2 minutes'),
    @(329, 328, 'Sunday, Jan 15', '8:50 AM', 'LO3832', 'Warsaw', '(WAW)', 'LOT ', 'E195', '(SP-LNO)', '8:58 AM', '0 hours, 8 minutes'),
    @(330, 329, 'Sunday, Jan 15', '9:50 AM', 'FR6844', 'Copenhagen', '(CPH)', 'Ryanair ', 'B738', '(SP-RKM)', '9:56 AM', '0 hours, 6 minutes'),
    @(331, 330, 'Sunday, Jan 15', '10:10 AM', 'DY1031', 'Bergen', '(BGO)', 'Norwegian ', 'B738', '(LN-DYM)', '10:12 AM', '0 hours, 2 minutes'),
    @(332, 331, 'Sunday, Jan 15', '11:10 AM', 'FR6109', 'Manchester', '(MAN)', 'Ryanair ', 'B738', '(SP-RSW)', '11:08 AM', '0 hours, -2 minutes'),
    @(333, 332, 'Sunday, Jan 15', '11:55 AM', 'LO3836', 'Warsaw', '(WAW)', 'LOT ', 'E190', '(SP-LMB)', '11:58 AM', '0 hours, 3 minutes'),
    @(334, 333, 'Sunday, Jan 15', '12:00 PM', 'W61675', 'Hamburg', '(HAM)', 'Wizz Air ', 'A21N', '(9H-WBU)', '12:09 PM', '0 hours, 9 minutes'),
    @(335, 334, 'Sunday, Jan 15', '12:25 PM', 'W61765', 'Alesund', '(AES)', 'Wizz Air ', 'A320', '(HA-LYS)', '12:26 PM', '0 hours, 1 minutes'),
    @(336, 335, 'Sunday, Jan 15', '12:30 PM', 'FR8510', 'Oslo', '(TRF)', 'Ryanair ', 'B738', '(SP-RSO)', '12:28 PM', '0 hours, -2 minutes'),
    @(337, 336, 'Sunday, Jan 15', '1:10 PM', 'LH1643', 'Munich', '(MUC)', 'Lufthansa ', 'CRJ9', '(D-ACNI)', '1:10 PM', '0 hours, 0 minutes'),
    @(338, 337, 'Sunday, Jan 15', '1:30 PM', 'FR6123', 'Edinburgh', '(EDI)', 'Ryanair ', 'B738', '(SP-RKM)', '1:40 PM', '0 hours, 10 minutes'),
    @(339, 338, 'Sunday, Jan 15', '1:45 PM', 'FR6103', 'Newcastle', '(NCL)', 'Ryanair ', 'B738', '(SP-RKQ)', '1:47 PM', '0 hours, 2 minutes'),
    @(340, 339, 'Sunday, Jan 15', '2:15 PM', 'LH1377', 'Frankfurt', '(FRA)', 'Lufthansa ', 'CRJ9', '(D-ACKI)', '2:37 PM', '0 hours, 22 minutes'),
    @(341, 340, 'Sunday, Jan 15', '2:50 PM', 'LO3838', 'Warsaw', '(WAW)', 'LOT ', 'E75S', '(SP-LIB)', '2:54 PM', '0 hours, 4 minutes'),
    @(342, 341, 'Sunday, Jan 15', '2:50 PM', 'W61747', 'Stavanger', '(SVG)', 'Wizz Air ', 'A320', '(HA-LYO)', '3:02 PM', '0 hours, 12 minutes'),
    @(343, 342, 'Sunday, Jan 15', '3:25 PM', 'FR2595', 'Valencia', '(VLC)', 'Ryanair ', 'B738', '(SP-RSL)', '3:28 PM', '0 hours, 3 minutes'),
    @(344, 343, 'Sunday, Jan 15', '3:35 PM', 'DY1055', 'Oslo', '(OSL)', 'Norwegian ', 'B738', '(LN-DYU)', '3:54 PM', '0 hours, 19 minutes'),
    @(345, 344, 'Sunday, Jan 15', '3:40 PM', 'W61751', 'Turku', '(TKU)', 'Wizz Air ', 'A320', '(HA-LYH)', '3:52 PM', '0 hours, 12 minutes'),
    @(346, 345, 'Sunday, Jan 15', '3:45 PM', 'SK760', 'Copenhagen', '(CPH)', 'SAS ', 'A20N', '(EI-SIA)', '3:56 PM', '0 hours, 11 minutes')
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value  = $r[1]   # A - NUMBER
    $ws.Cells.Item($row, 2).Value  = $r[2]   # B - DATE
    $ws.Cells.Item($row, 3).Value  = $r[3]   # C - TIME
    $ws.Cells.Item($row, 4).Value  = $r[4]   # D - FLIGHT
    $ws.Cells.Item($row, 5).Value  = $r[5]   # E - TO
    $ws.Cells.Item($row, 6).Value  = $r[6]   # F - SHORT
    $ws.Cells.Item($row, 7).Value  = $r[7]   # G - AIRLINE
    $ws.Cells.Item($row, 8).Value  = $r[8]   # H - MODEL
    $ws.Cells.Item($row, 9).Value  = $r[9]   # I - AIRCFAT ID
    $ws.Cells.Item($row, 10).Value = $r[10]  # J - STATUS
    $ws.Cells.Item($row, 12).Value = $r[11]  # L - DIFFERENCE
}
